$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# --- Crime Complaints table updates (rows 15-30) ---
# Row 15
$ws.Range("C23").Copy($ws.Range("D15"))
$ws.Range("E23").Copy($ws.Range("E15"))
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -33.333333333333
$ws.Range("L15").Value = -11.111111111111
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -27.272727272727

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 250
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 31.25
$ws.Range("I16").Value = 168
$ws.Range("J16").Value = 220
$ws.Range("K16").Value = -23.636363636363
$ws.Range("L16").Value = 3.067484662576
$ws.Range("M16").Value = 28.24427480916
$ws.Range("N16").Value = -77.689243027888

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -28.571428571428
$ws.Range("I17").Value = 176
$ws.Range("J17").Value = 165
$ws.Range("K17").Value = 6.666666666666
$ws.Range("L17").Value = 10
$ws.Range("M17").Value = 109.52380952381
$ws.Range("N17").Value = -33.584905660377

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = -34.285714285714
$ws.Range("I18").Value = 251
$ws.Range("J18").Value = 385
$ws.Range("K18").Value = -34.805194805194
$ws.Range("L18").Value = 18.957345971564
$ws.Range("M18").Value = 64.052287581699
$ws.Range("N18").Value = -65.187239944521

# Row 19
$ws.Range("C19").Value = 33
$ws.Range("D19").Value = 30
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 156
$ws.Range("G19").Value = 158
$ws.Range("H19").Value = -1.26582278481
$ws.Range("I19").Value = 1167
$ws.Range("J19").Value = 1202
$ws.Range("K19").Value = -2.911813643926
$ws.Range("L19").Value = 60.522696011004
$ws.Range("M19").Value = 28.241758241758
$ws.Range("N19").Value = -48.40848806366

# Row 20
$ws.Range("C23").Copy($ws.Range("C20"))
$ws.Range("E20").Value = -100
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = -20.37037037037
$ws.Range("M20").Value = 7.5
$ws.Range("N20").Value = -93.610698365527

# Row 21
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 214
$ws.Range("G21").Value = 227
$ws.Range("H21").Value = -5.726872246696
$ws.Range("I21").Value = 1814
$ws.Range("J21").Value = 2040
$ws.Range("K21").Value = -11.078431372549
$ws.Range("L21").Value = 37.946768060836
$ws.Range("M21").Value = 36.390977443609
$ws.Range("N21").Value = -61.288945795988

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 39
$ws.Range("J22").Value = 37
$ws.Range("K22").Value = 5.405405405405
$ws.Range("L22").Value = 44.444444444444
$ws.Range("M22").Value = -9.302325581395

# Row 24
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 173
$ws.Range("G24").Value = 135
$ws.Range("H24").Value = 28.148148148148
$ws.Range("I24").Value = 1835
$ws.Range("J24").Value = 1706
$ws.Range("K24").Value = 7.561547479484
$ws.Range("L24").Value = 63.111111111111
$ws.Range("M24").Value = 39.437689969604

# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -27.272727272727
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = -31.25
$ws.Range("I25").Value = 376
$ws.Range("J25").Value = 368
$ws.Range("K25").Value = 2.173913043478
$ws.Range("L25").Value = 37.226277372262
$ws.Range("M25").Value = 62.068965517241

# Row 26
$ws.Range("C23").Copy($ws.Range("D26"))
$ws.Range("E23").Copy($ws.Range("E26"))

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C23").Copy($ws.Range("D27"))
$ws.Range("E23").Copy($ws.Range("E27"))
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 58
$ws.Range("K27").Value = -13.432835820895
$ws.Range("L27").Value = -4.918032786885

# Row 30
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 12
$ws.Range("K30").Value = -33.333333333333
